$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data for columns I ("I0") and J ("IF"), rows 2-45
$dataI0 = @(5,5,10,6,8,1,9,8,8,8,5,8,6,7,9,6,9,5,6,4,7,5,6,7,7,8,8,9,9,9,8,9,7,5,8,8,8,9,8,8,8,4,7,2)
$dataIF = @(6,6,10,7,8,1,9,8,8,8,5,8,7,8,9,6,9,5,7,5,7,6,7,8,7,9,8,9,9,9,9,9,7,6,8,9,8,9,8,9,8,5,7,2)

for ($i = 0; $i -lt $dataI0.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI0[$i]
    $ws.Cells.Item($row, 10).Value = $dataIF[$i]
}

$null = $ws.Range("A1").Select()
